# Auto-generated edit script: apply updated market-price / profit values
# to the Kujata_Profits leve-profit sheets (scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 185.14285
$ws.Range("I38").Value = 185.14285
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 555.4285500000001
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -183.4285500000001
$ws.Range("N38").Value = ""

$ws.Range("H58").Value = 1131.2
$ws.Range("I58").Value = 417.25
$ws.Range("J58").Value = 2202.125
$ws.Range("K58").Value = 1251.75
$ws.Range("L58").Value = 6606.375
$ws.Range("M58").Value = -1101.75
$ws.Range("N58").Value = -6906.375

$ws.Range("H64").Value = 3710
$ws.Range("I64").Value = 3663.3333
$ws.Range("J64").Value = 3780
$ws.Range("K64").Value = 3663.3333
$ws.Range("L64").Value = 3780
$ws.Range("M64").Value = -3415.3333
$ws.Range("N64").Value = -4276

$ws.Range("H67").Value = 3710
$ws.Range("I67").Value = 3663.3333
$ws.Range("J67").Value = 3780
$ws.Range("K67").Value = 3663.3333
$ws.Range("L67").Value = 3780
$ws.Range("M67").Value = -2805.3333
$ws.Range("N67").Value = -5496

$ws.Range("H100").Value = 1487.5
$ws.Range("I100").Value = 1100
$ws.Range("J100").Value = 2650
$ws.Range("K100").Value = 1100
$ws.Range("L100").Value = 2650
$ws.Range("M100").Value = -559
$ws.Range("N100").Value = -3732

$ws.Range("H137").Value = 1465.6
$ws.Range("I137").Value = 1188.4667
$ws.Range("J137").Value = 1604.1666
$ws.Range("K137").Value = 3565.4001
$ws.Range("L137").Value = 4812.4998
$ws.Range("M137").Value = -1015.4001
$ws.Range("N137").Value = -9912.4998

$ws.Range("H138").Value = 445539.44
$ws.Range("I138").Value = 1397.3334
$ws.Range("J138").Value = 525806.0600000001
$ws.Range("K138").Value = 4192.0002
$ws.Range("L138").Value = 1577418.18
$ws.Range("M138").Value = 947.9997999999996
$ws.Range("N138").Value = -1587698.18

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5553.0757
$ws.Range("I32").Value = 5417.9214
$ws.Range("J32").Value = 8999.5
$ws.Range("K32").Value = 5417.9214
$ws.Range("L32").Value = 8999.5
$ws.Range("M32").Value = -5130.9214
$ws.Range("N32").Value = -9573.5

$ws.Range("H63").Value = 1863.1052
$ws.Range("I63").Value = 1665.4482
$ws.Range("J63").Value = 2500
$ws.Range("K63").Value = 1665.4482
$ws.Range("L63").Value = 2500
$ws.Range("M63").Value = -979.4482
$ws.Range("N63").Value = -3872

$ws.Range("H66").Value = 1863.1052
$ws.Range("I66").Value = 1665.4482
$ws.Range("J66").Value = 2500
$ws.Range("K66").Value = 8327.241
$ws.Range("L66").Value = 12500
$ws.Range("M66").Value = -4895.241
$ws.Range("N66").Value = -19364

$ws.Range("H74").Value = 2188.9167
$ws.Range("I74").Value = 1426.7
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 1426.7
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -552.7
$ws.Range("N74").Value = -7748

$ws.Range("H77").Value = 2188.9167
$ws.Range("I77").Value = 1426.7
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 7133.5
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -2765.5
$ws.Range("N77").Value = -38736

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 623.875
$ws.Range("I80").Value = 300
$ws.Range("J80").Value = 698.61536
$ws.Range("K80").Value = 300
$ws.Range("L80").Value = 698.61536
$ws.Range("M80").Value = 698
$ws.Range("N80").Value = -2694.61536

$ws.Range("H83").Value = 623.875
$ws.Range("I83").Value = 300
$ws.Range("J83").Value = 698.61536
$ws.Range("K83").Value = 1500
$ws.Range("L83").Value = 3493.0768
$ws.Range("M83").Value = 3492
$ws.Range("N83").Value = -13477.0768

$ws.Range("H105").Value = 48091220
$ws.Range("I105").Value = 50495684
$ws.Range("J105").Value = 1911
$ws.Range("K105").Value = 50495684
$ws.Range("L105").Value = 1911
$ws.Range("M105").Value = -50493937
$ws.Range("N105").Value = -5405

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 350351
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 350351
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 350351
$ws.Range("M22").Value = ""
$ws.Range("N22").Value = -351051

$ws.Range("H31").Value = 1645.7028
$ws.Range("I31").Value = 1443.9688
$ws.Range("J31").Value = 2936.8
$ws.Range("K31").Value = 1443.9688
$ws.Range("L31").Value = 2936.8
$ws.Range("M31").Value = -1148.9688
$ws.Range("N31").Value = -3526.8

$ws.Range("H32").Value = 5000
$ws.Range("I32").Value = 5000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 5000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4684

$ws.Range("H34").Value = 1645.7028
$ws.Range("I34").Value = 1443.9688
$ws.Range("J34").Value = 2936.8
$ws.Range("K34").Value = 1443.9688
$ws.Range("L34").Value = 2936.8
$ws.Range("M34").Value = -1241.9688
$ws.Range("N34").Value = -3340.8

$ws.Range("H35").Value = 275
$ws.Range("I35").Value = 212.5
$ws.Range("J35").Value = 400
$ws.Range("K35").Value = 212.5
$ws.Range("L35").Value = 400
$ws.Range("M35").Value = 81.5
$ws.Range("N35").Value = -988

$ws.Range("H38").Value = 2750
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 2750
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 2750
$ws.Range("N38").Value = -3504

$ws.Range("H42").Value = 2000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 2000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 2000
$ws.Range("N42").Value = -3186

$ws.Range("H45").Value = 3900
$ws.Range("I45").Value = 3900
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 3900
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -3307

$ws.Range("H46").Value = 2750
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2750
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2750
$ws.Range("N46").Value = -3172

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 5313.95
$ws.Range("I107").Value = 409.16666
$ws.Range("J107").Value = 7416
$ws.Range("K107").Value = 1227.49998
$ws.Range("L107").Value = 22248
$ws.Range("M107").Value = 692.5000199999999
$ws.Range("N107").Value = -26088

$ws.Range("H122").Value = 1659.6842
$ws.Range("I122").Value = 850
$ws.Range("J122").Value = 1704.6666
$ws.Range("K122").Value = 7650
$ws.Range("L122").Value = 15341.9994
$ws.Range("M122").Value = -5200
$ws.Range("N122").Value = -20241.9994

$ws.Range("H131").Value = 14286622
$ws.Range("I131").Value = 142857550
$ws.Range("J131").Value = 963.20636
$ws.Range("K131").Value = 428572650
$ws.Range("L131").Value = 2889.61908
$ws.Range("M131").Value = -428567610
$ws.Range("N131").Value = -12969.61908

$ws.Range("H139").Value = 1811.2162
$ws.Range("I139").Value = 1879.2174
$ws.Range("J139").Value = 1699.5
$ws.Range("K139").Value = 5637.6522
$ws.Range("L139").Value = 5098.5
$ws.Range("M139").Value = -497.6522000000004
$ws.Range("N139").Value = -15378.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 25003676
$ws.Range("I70").Value = 25003808
$ws.Range("J70").Value = 25003512
$ws.Range("K70").Value = 25003808
$ws.Range("L70").Value = 25003512
$ws.Range("M70").Value = -25003538
$ws.Range("N70").Value = -25004052

$ws.Range("H73").Value = 25003676
$ws.Range("I73").Value = 25003808
$ws.Range("J73").Value = 25003512
$ws.Range("K73").Value = 25003808
$ws.Range("L73").Value = 25003512
$ws.Range("M73").Value = -25002872
$ws.Range("N73").Value = -25005384

$ws.Range("H135").Value = 36133.332
$ws.Range("I135").Value = 50000
$ws.Range("J135").Value = 34400
$ws.Range("K135").Value = 50000
$ws.Range("L135").Value = 34400
$ws.Range("M135").Value = -44930
$ws.Range("N135").Value = -44540

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 4125
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 4125
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 4125
$ws.Range("N94").Value = -5477

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 27779664
$ws.Range("I122").Value = 27779664
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 83338992
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -83336542

$ws.Range("H136").Value = 1339.1154
$ws.Range("I136").Value = 1294.8823
$ws.Range("J136").Value = 1422.6666
$ws.Range("K136").Value = 3884.6469
$ws.Range("L136").Value = 4267.9998
$ws.Range("M136").Value = 347.8696199999999
$ws.Range("N136").Value = -1334.6469
